$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '63.551.18'
$ws.Cells.Item(2, 5).Value = '  -3.39%  '
# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.603.81'
$ws.Cells.Item(3, 5).Value = '  -2.22%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '573.84'
$ws.Cells.Item(5, 5).Value = '  -4.10%  '
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '154.77'
$ws.Cells.Item(6, 5).Value = '  -1.65%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.06%  '
# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.621'
$ws.Cells.Item(8, 5).Value = '  -5.02%  '
# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.117'
$ws.Cells.Item(9, 5).Value = '  -6.81%  '
# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.46%  '
# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.381'
$ws.Cells.Item(11, 5).Value = '  -5.23%  '
# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.60%  '
# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '28.20'
$ws.Cells.Item(13, 5).Value = '  -2.17%  '
# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '3.073.86'
$ws.Cells.Item(14, 5).Value = '  -1.92%  '
# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.0000177'
$ws.Cells.Item(15, 5).Value = '  -7.93%  '
# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '63.366.74'
$ws.Cells.Item(16, 5).Value = '  -3.52%  '
# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '2.606.26'
$ws.Cells.Item(17, 5).Value = '  -1.33%  '
# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '11.96'
$ws.Cells.Item(18, 5).Value = '  -4.86%  '
# Row 19
$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.53'
$ws.Cells.Item(19, 5).Value = '  -5.26%  '
# Row 20
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.46'
$ws.Cells.Item(20, 5).Value = '  +0.64%  '
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '342.30'
$ws.Cells.Item(21, 5).Value = '  -1.77%  '
# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.05%  '
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '66.95'
$ws.Cells.Item(23, 5).Value = '  -3.82%  '
# Row 24
$ws.Cells.Item(24, 5).Value = '  -3.77%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  -3.97%  '
# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '586.41'
$ws.Cells.Item(26, 5).Value = '  +2.82%  '
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.11'
$ws.Cells.Item(27, 5).Value = '  -4.66%  '
# Row 28
$ws.Cells.Item(28, 5).Value = '  -4.17%  '
# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.01%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.82%  '
# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.87'
$ws.Cells.Item(31, 5).Value = '  -2.74%  '
# Row 32
$ws.Cells.Item(32, 5).Value = '  -4.51%  '
# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.72'
$ws.Cells.Item(33, 5).Value = '  -5.31%  '
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '6.51'
$ws.Cells.Item(34, 5).Value = '  -2.71%  '
# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.36'
$ws.Cells.Item(35, 5).Value = '  -1.58%  '
# Row 36
$ws.Cells.Item(36, 5).Value = '  -4.24%  '
# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.998'
$ws.Cells.Item(37, 5).Value = '  -0.05%  '
# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '19.62'
$ws.Cells.Item(38, 5).Value = '  -4.58%  '
# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '154.78'
$ws.Cells.Item(39, 5).Value = '  -0.38%  '
# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.85'
$ws.Cells.Item(40, 5).Value = '  -4.08%  '
# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.999'
$ws.Cells.Item(41, 5).Value = '  -0.03%  '
# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '41.41'
# Row 43
$ws.Cells.Item(43, 5).Value = '  +6.78%  '
# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '155.52'
# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.88'
$ws.Cells.Item(45, 5).Value = '  -4.88%  '
# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '23.09'
$ws.Cells.Item(46, 5).Value = '  +1.40%  '
# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0585'
$ws.Cells.Item(47, 5).Value = '  -4.40%  '
# Row 48
$ws.Cells.Item(48, 5).Value = '  -1.97%  '
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.100'
$ws.Cells.Item(49, 5).Value = '  -1.65%  '
# Row 50
$ws.Cells.Item(50, 5).Value = '  -3.78%  '
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '18.77'
$ws.Cells.Item(51, 5).Value = '  -5.22%  '
